$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Title heading (appears twice: main heading and bold summary line near the end)
Replace-Text "Play Hitman Slot for Free - Exciting Bonus Games Included" "Play Hitman Free: Exciting Bonus Games & Attractive RTP"

# "What we like" bullet list changes
Replace-Text "Multiple winning combinations" "Engaging gameplay"
Replace-Text "Straightforward gameplay" "Good payout average"

# "What we don't like" bullet list changes
Replace-Text "Graphics could be improved" "Graphics could be better"
Replace-Text "No additional free spins during free spin phase" "Limited free spins"

# Final italic summary paragraph
Replace-Text "Read our review of Hitman, a Microgaming-made online slot game with 15 paylines and bonus features. Play it for free and enjoy exciting bonus games." "Play Hitman for free and enjoy its exciting bonus games and attractive RTP value."
